$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (preserve trailing zeros / formatting)
$ws.Range("D2").Value = '67.321.88'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '3.512.31'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.59'
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.93'
$ws.Range("E6").Value = '  +2.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.32%  '
$ws.Range("E9").Value = '  +2.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.18'
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.432'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '4.112.34'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.13'
$ws.Range("E14").Value = '  +6.98%  '
$ws.Range("D15").Value = '67.285.12'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '3.504.63'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.30'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.50'
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '393.49'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.60'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.536'
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.69'
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.16'
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.40'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.46'
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.879'
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.93'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.90'
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.24'
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0732'
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.814.53'
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.12'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.56'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.56'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0304'
$ws.Range("E47").Value = '  -2.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '343.13'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  -0.77%  '
